$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "2020" data column (H) -------------------------------
# Each row's new H cell must inherit the same cell formatting (number
# format / borders / font) as its existing neighbours. Copy the format
# from the most appropriate existing cell, then (for non-blank rows)
# write the value.

function Set-FormatFrom($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

# row : (value or $null for a blank cell, source cell to copy formatting from)
Set-FormatFrom "G4" "H4"
$ws.Range("H4").Value = 2020

Set-FormatFrom "G5" "H5"
$ws.Range("H5").Value = 42.2

Set-FormatFrom "G6" "H6"

Set-FormatFrom "G7" "H7"
$ws.Range("H7").Value = 42.5

Set-FormatFrom "G10" "H8"
$ws.Range("H8").Value = 42

Set-FormatFrom "G9" "H9"

Set-FormatFrom "G10" "H10"
$ws.Range("H10").Value = 50.9

Set-FormatFrom "G11" "H11"
$ws.Range("H11").Value = 36.9

Set-FormatFrom "G12" "H12"
$ws.Range("H12").Value = 34.8

Set-FormatFrom "G13" "H13"

Set-FormatFrom "G14" "H14"
$ws.Range("H14").Value = 30.7

Set-FormatFrom "G15" "H15"
$ws.Range("H15").Value = 48.8

Set-FormatFrom "G16" "H16"

Set-FormatFrom "G17" "H17"
$ws.Range("H17").Value = 61.1

Set-FormatFrom "G18" "H18"
$ws.Range("H18").Value = 56.7

Set-FormatFrom "G19" "H19"
$ws.Range("H19").Value = 41.6

Set-FormatFrom "G22" "H20"
$ws.Range("H20").Value = 49

Set-FormatFrom "G21" "H21"
$ws.Range("H21").Value = 43.5

Set-FormatFrom "G22" "H22"
$ws.Range("H22").Value = 33.9

Set-FormatFrom "G23" "H23"
$ws.Range("H23").Value = 34.6

Set-FormatFrom "G24" "H24"
$ws.Range("H24").Value = 23.6

Set-FormatFrom "G25" "H25"
$ws.Range("H25").Value = 35.9

# --- Update the active view / selection --------------------------------
# Original view had topLeftCell="B1" and selection H15; target view has
# no explicit topLeftCell (so it defaults back to A1) and selection B13.
$ws.Range("B13").Select() | Out-Null

Write-Host "Edit applied"
